$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A3's date style: it currently uses a custom "d-mmm" style (numFmtId 16).
# Replace it with the same short-date style already used by A2/A4.. (numFmtId 14)
# by copying A2's format onto A3 (reuses the existing style index instead of
# allocating a new custom number format).
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Append the new time-tracking entry as row 11 ---
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = "accountmanage seit mit passwort ändern,reservierungen löschbar, details,tauschseite und anfragenseite(noch ohne datenbankanbindung)"

# Give A11 the same date style as the other date cells, then set its value
# (serial date number for 2017-01-09).
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("A11").Value = 42731

# --- Update the selected cell shown in the sheet view ---
[void]$ws.Range("C11").Select()
